$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Level" column (C) values for a few rows
$ws.Range("C2").Value = "moderate"
$ws.Range("C3").Value = "moderate"
$ws.Range("C5").Value = "easy"
